# Scheduled-runner style refresh of market/profit figures across the
# "Leve Profit" tables on each job sheet (ALC, ARM, CRP, CUL, GSM, LTW, WVR).
# For every touched row, columns H:N (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) are
# refreshed with new market-derived values. Where a column's new value
# would be an empty/absent cell (as in the source data refresh), the
# cell content is cleared instead of being set to 0.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 93.75
$ws.Range("I2").Value = 90
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 90
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 23
$ws.Range("N2").Value = -326
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10460
$ws.Range("M20").ClearContents()
$ws.Range("H35").Value = 10000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 10000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 10000
$ws.Range("N35").Value = -10758
$ws.Range("M35").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H112").Value = 1741.225
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1806.5526
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 5419.6578
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -7635.6578
$ws.Range("H129").Value = 899.4
$ws.Range("I129").Value = 495.7143
$ws.Range("J129").Value = 965.1163
$ws.Range("K129").Value = 1487.1429
$ws.Range("L129").Value = 2895.3489
$ws.Range("M129").Value = 3512.8571
$ws.Range("H138").Value = 3150603.8
$ws.Range("I138").Value = 12000
$ws.Range("J138").Value = 3210961.5
$ws.Range("K138").Value = 36000
$ws.Range("L138").Value = 9632884.5
$ws.Range("M138").Value = -30860
$ws.Range("N138").Value = -9643164.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3060.3462
$ws.Range("I45").Value = 3181.6191
$ws.Range("J45").Value = 2551
$ws.Range("K45").Value = 3181.6191
$ws.Range("L45").Value = 2551
$ws.Range("M45").Value = -2804.6191
$ws.Range("N45").Value = -3305
$ws.Range("H61").Value = 37112080
$ws.Range("I61").Value = 50051012
$ws.Range("J61").Value = 143704
$ws.Range("K61").Value = 50051012
$ws.Range("L61").Value = 143704
$ws.Range("M61").Value = -50050800
$ws.Range("N61").Value = -144128
$ws.Range("H74").Value = 7201054.5
$ws.Range("I74").Value = 14765422
$ws.Range("J74").Value = 56930
$ws.Range("K74").Value = 14765422
$ws.Range("L74").Value = 56930
$ws.Range("M74").Value = -14764548
$ws.Range("N74").Value = -58678
$ws.Range("H77").Value = 7201054.5
$ws.Range("I77").Value = 14765422
$ws.Range("J77").Value = 56930
$ws.Range("K77").Value = 73827110
$ws.Range("L77").Value = 284650
$ws.Range("M77").Value = -73822742
$ws.Range("N77").Value = -293386
$ws.Range("H97").Value = 2978038.5
$ws.Range("I97").Value = 3290884.8
$ws.Range("J97").Value = 6000
$ws.Range("K97").Value = 3290884.8
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -3290388.8
$ws.Range("N97").Value = -6992
$ws.Range("H132").Value = 68663.164
$ws.Range("I132").Value = 51545.75
$ws.Range("J132").Value = 102898
$ws.Range("K132").Value = 154637.25
$ws.Range("L132").Value = 308694
$ws.Range("M132").Value = -152107.25
$ws.Range("N132").Value = -313754
$ws.Range("H136").Value = 37112080
$ws.Range("I136").Value = 50051012
$ws.Range("J136").Value = 143704
$ws.Range("K136").Value = 150153036
$ws.Range("L136").Value = 431112
$ws.Range("M136").Value = -150150486
$ws.Range("N136").Value = -436212

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 46833.332
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 46833.332
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 46833.332
$ws.Range("N52").Value = -47421.332
$ws.Range("H122").Value = 1315.8529
$ws.Range("I122").Value = 1148.0938
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3444.2814
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -994.2814000000003
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H68").Value = 329249.56
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 329249.56
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 987748.6799999999
$ws.Range("N68").Value = -989370.6799999999
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 329249.56
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 329249.56
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 2963246.04
$ws.Range("N71").Value = -2971358.04
$ws.Range("M71").ClearContents()
$ws.Range("H107").Value = 1070.6742
$ws.Range("I107").Value = 504.5
$ws.Range("J107").Value = 1532.8572
$ws.Range("K107").Value = 1513.5
$ws.Range("L107").Value = 4598.571599999999
$ws.Range("M107").Value = 406.5
$ws.Range("N107").Value = -8438.571599999999
$ws.Range("H113").Value = 568.225
$ws.Range("I113").Value = 486.53845
$ws.Range("J113").Value = 607.55554
$ws.Range("K113").Value = 1459.61535
$ws.Range("L113").Value = 1822.66662
$ws.Range("M113").Value = 710.38465
$ws.Range("H117").Value = 3031323
$ws.Range("I117").Value = 576.3333
$ws.Range("J117").Value = 4167853
$ws.Range("K117").Value = 1728.9999
$ws.Range("L117").Value = 12503559
$ws.Range("M117").Value = 1713.0001
$ws.Range("N117").Value = -12510443
$ws.Range("H122").Value = 529.375
$ws.Range("I122").Value = 355.14285
$ws.Range("J122").Value = 1749
$ws.Range("K122").Value = 3196.28565
$ws.Range("L122").Value = 15741
$ws.Range("M122").Value = -746.2856500000003
$ws.Range("N122").Value = -20641
$ws.Range("H129").Value = 2382657.2
$ws.Range("I129").Value = 785.6087
$ws.Range("J129").Value = 6947911.5
$ws.Range("K129").Value = 2356.8261
$ws.Range("L129").Value = 20843734.5
$ws.Range("M129").Value = 2643.1739
$ws.Range("N129").Value = -20853734.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 6609.2856
$ws.Range("I99").Value = 4810
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 4810
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -2564
$ws.Range("N99").Value = -34492
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 456.47058
$ws.Range("I22").Value = 241.6
$ws.Range("J22").Value = 763.4286
$ws.Range("K22").Value = 241.6
$ws.Range("L22").Value = 763.4286
$ws.Range("M22").Value = 53.40000000000001
$ws.Range("N22").Value = -1353.4286
$ws.Range("H27").Value = 456.47058
$ws.Range("I27").Value = 241.6
$ws.Range("J27").Value = 763.4286
$ws.Range("K27").Value = 241.6
$ws.Range("L27").Value = 763.4286
$ws.Range("M27").Value = -134.6
$ws.Range("N27").Value = -977.4286
$ws.Range("H74").Value = 3344732.2
$ws.Range("I74").Value = 5005098.5
$ws.Range("J74").Value = 24000
$ws.Range("K74").Value = 5005098.5
$ws.Range("L74").Value = 24000
$ws.Range("M74").Value = -5004100.5
$ws.Range("N74").Value = -25996
$ws.Range("H77").Value = 3344732.2
$ws.Range("I77").Value = 5005098.5
$ws.Range("J77").Value = 24000
$ws.Range("K77").Value = 15015295.5
$ws.Range("L77").Value = 72000
$ws.Range("M77").Value = -15010303.5
$ws.Range("N77").Value = -81984
$ws.Range("H93").Value = 2056.3
$ws.Range("I93").Value = 1816.6
$ws.Range("J93").Value = 2296
$ws.Range("K93").Value = 1816.6
$ws.Range("L93").Value = 2296
$ws.Range("M93").Value = -568.5999999999999
$ws.Range("N93").Value = -4792
$ws.Range("H132").Value = 115311.336
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 115311.336
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 345934.008
$ws.Range("N132").Value = -350994.008
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 109736.7
$ws.Range("I136").Value = 55607.21
$ws.Range("J136").Value = 238294.25
$ws.Range("K136").Value = 166821.63
$ws.Range("L136").Value = 714882.75
$ws.Range("M136").Value = -164271.63

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2616.9167
$ws.Range("I81").Value = 1566.3334
$ws.Range("J81").Value = 2967.111
$ws.Range("K81").Value = 3132.6668
$ws.Range("L81").Value = 5934.222
$ws.Range("M81").Value = -2071.6668
$ws.Range("N81").Value = -8056.222
$ws.Range("H84").Value = 2616.9167
$ws.Range("I84").Value = 1566.3334
$ws.Range("J84").Value = 2967.111
$ws.Range("K84").Value = 15663.334
$ws.Range("L84").Value = 29671.11
$ws.Range("M84").Value = -10359.334
$ws.Range("N84").Value = -40279.11
$ws.Range("H126").Value = 2287.375
$ws.Range("I126").Value = 1199.75
$ws.Range("J126").Value = 3375
$ws.Range("K126").Value = 3599.25
$ws.Range("L126").Value = 10125
$ws.Range("M126").Value = -1129.25
$ws.Range("N126").Value = -15065
$ws.Range("H132").Value = 170125.42
$ws.Range("I132").Value = 145944
$ws.Range("J132").Value = 203979.4
$ws.Range("K132").Value = 437832
$ws.Range("L132").Value = 611938.2
$ws.Range("M132").Value = -435302
$ws.Range("H136").Value = 296872.84
$ws.Range("I136").Value = 261000
$ws.Range("J136").Value = 344703.34
$ws.Range("K136").Value = 783000
$ws.Range("L136").Value = 1034110.02
$ws.Range("M136").Value = -780450
$ws.Range("N136").Value = -1039210.02
